# Efna4 -> Epha4 NATMI edge table: refreshed with new TPM-derived values.
# Sending/target clusters are now alphabetically ordered (ECs, FAPs, Inflammatory-Mac,
# MuSCs, Resolving-Mac); "Resolving-Mac" is a sending cluster and "Inflammatory-Mac"
# is now also a valid target cluster, growing the table from 16 to 20 data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 20,20

# Row 2: ECs -> ECs (Efna4/Epha4)
$data[0,0] = "ECs"
$data[0,1] = "Efna4"
$data[0,2] = "Epha4"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 1.974568
$data[0,7] = 5.923704
$data[0,8] = 0.5990695552080697
$data[0,9] = 0.5990695552080698
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 8.540560666666666
$data[0,13] = 25.621682
$data[0,14] = 0.4159358086620884
$data[0,15] = 0.4159358086620884
$data[0,16] = 16.86391779445866
$data[0,17] = 151.775260150128
$data[0,18] = 0.2491744798903061
$data[0,19] = 0.2491744798903061

# Row 3: ECs -> FAPs (Efna4/Epha4)
$data[1,0] = "ECs"
$data[1,1] = "Efna4"
$data[1,2] = "Epha4"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 1.974568
$data[1,7] = 5.923704
$data[1,8] = 0.5990695552080697
$data[1,9] = 0.5990695552080698
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 11.57455166666667
$data[1,13] = 34.723655
$data[1,14] = 0.563694901924408
$data[1,15] = 0.563694901924408
$data[1,16] = 22.85473933534666
$data[1,17] = 205.69265401812
$data[1,18] = 0.3376924541689116
$data[1,19] = 0.3376924541689116

# Row 4: ECs -> Inflammatory-Mac (Efna4/Epha4)
$data[2,0] = "ECs"
$data[2,1] = "Efna4"
$data[2,2] = "Epha4"
$data[2,3] = "Inflammatory-Mac"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 1.974568
$data[2,7] = 5.923704
$data[2,8] = 0.5990695552080697
$data[2,9] = 0.5990695552080698
$data[2,10] = 1
$data[2,11] = 0.3333333333333333
$data[2,12] = 0.026642
$data[2,13] = 0.079926
$data[2,14] = 0.001297498167494471
$data[2,15] = 0.001297498167494471
$data[2,16] = 0.05260644065599999
$data[2,17] = 0.473457965904
$data[2,18] = 0.0007772916500841984
$data[2,19] = 0.0007772916500841985

# Row 5: ECs -> MuSCs (Efna4/Epha4)
$data[3,0] = "ECs"
$data[3,1] = "Efna4"
$data[3,2] = "Epha4"
$data[3,3] = "MuSCs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 1.974568
$data[3,7] = 5.923704
$data[3,8] = 0.5990695552080697
$data[3,9] = 0.5990695552080698
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 0.391608
$data[3,13] = 1.174824
$data[3,14] = 0.01907179124600912
$data[3,15] = 0.01907179124600912
$data[3,16] = 0.773256625344
$data[3,17] = 6.959309628096
$data[3,18] = 0.01142532949876784
$data[3,19] = 0.01142532949876784

# Row 6: FAPs -> ECs (Efna4/Epha4)
$data[4,0] = "FAPs"
$data[4,1] = "Efna4"
$data[4,2] = "Epha4"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 1.020259
$data[4,7] = 3.060777
$data[4,8] = 0.309539152527049
$data[4,9] = 0.309539152527049
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 8.540560666666666
$data[4,13] = 25.621682
$data[4,14] = 0.4159358086620884
$data[4,15] = 0.4159358086620884
$data[4,16] = 8.713583885212666
$data[4,17] = 78.422254966914
$data[4,18] = 0.1287484177189157
$data[4,19] = 0.1287484177189157

# Row 7: FAPs -> FAPs (Efna4/Epha4)
$data[5,0] = "FAPs"
$data[5,1] = "Efna4"
$data[5,2] = "Epha4"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 1.020259
$data[5,7] = 3.060777
$data[5,8] = 0.309539152527049
$data[5,9] = 0.309539152527049
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 11.57455166666667
$data[5,13] = 34.723655
$data[5,14] = 0.563694901924408
$data[5,15] = 0.563694901924408
$data[5,16] = 11.80904050888167
$data[5,17] = 106.281364579935
$data[5,18] = 0.1744856422254993
$data[5,19] = 0.1744856422254993

# Row 8: FAPs -> Inflammatory-Mac (Efna4/Epha4)
$data[6,0] = "FAPs"
$data[6,1] = "Efna4"
$data[6,2] = "Epha4"
$data[6,3] = "Inflammatory-Mac"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 1.020259
$data[6,7] = 3.060777
$data[6,8] = 0.309539152527049
$data[6,9] = 0.309539152527049
$data[6,10] = 1
$data[6,11] = 0.3333333333333333
$data[6,12] = 0.026642
$data[6,13] = 0.079926
$data[6,14] = 0.001297498167494471
$data[6,15] = 0.001297498167494471
$data[6,16] = 0.027181740278
$data[6,17] = 0.244635662502
$data[6,18] = 0.0004016264831716377
$data[6,19] = 0.0004016264831716377

# Row 9: FAPs -> MuSCs (Efna4/Epha4)
$data[7,0] = "FAPs"
$data[7,1] = "Efna4"
$data[7,2] = "Epha4"
$data[7,3] = "MuSCs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 1.020259
$data[7,7] = 3.060777
$data[7,8] = 0.309539152527049
$data[7,9] = 0.309539152527049
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 0.391608
$data[7,13] = 1.174824
$data[7,14] = 0.01907179124600912
$data[7,15] = 0.01907179124600912
$data[7,16] = 0.399541586472
$data[7,17] = 3.595874278248
$data[7,18] = 0.005903466099462454
$data[7,19] = 0.005903466099462454

# Row 10: Inflammatory-Mac -> ECs (Efna4/Epha4)
$data[8,0] = "Inflammatory-Mac"
$data[8,1] = "Efna4"
$data[8,2] = "Epha4"
$data[8,3] = "ECs"
$data[8,4] = 1
$data[8,5] = 0.3333333333333333
$data[8,6] = 0.077601
$data[8,7] = 0.232803
$data[8,8] = 0.02354357842003994
$data[8,9] = 0.02354357842003994
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 8.540560666666666
$data[8,13] = 25.621682
$data[8,14] = 0.4159358086620884
$data[8,15] = 0.4159358086620884
$data[8,16] = 0.662756048294
$data[8,17] = 5.964804434646
$data[8,18] = 0.009792617328938606
$data[8,19] = 0.009792617328938607

# Row 11: Inflammatory-Mac -> FAPs (Efna4/Epha4)
$data[9,0] = "Inflammatory-Mac"
$data[9,1] = "Efna4"
$data[9,2] = "Epha4"
$data[9,3] = "FAPs"
$data[9,4] = 1
$data[9,5] = 0.3333333333333333
$data[9,6] = 0.077601
$data[9,7] = 0.232803
$data[9,8] = 0.02354357842003994
$data[9,9] = 0.02354357842003994
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 11.57455166666667
$data[9,13] = 34.723655
$data[9,14] = 0.563694901924408
$data[9,15] = 0.563694901924408
$data[9,16] = 0.898196783885
$data[9,17] = 8.083771054965
$data[9,18] = 0.01327139512843402
$data[9,19] = 0.01327139512843402

# Row 12: Inflammatory-Mac -> Inflammatory-Mac (Efna4/Epha4)
$data[10,0] = "Inflammatory-Mac"
$data[10,1] = "Efna4"
$data[10,2] = "Epha4"
$data[10,3] = "Inflammatory-Mac"
$data[10,4] = 1
$data[10,5] = 0.3333333333333333
$data[10,6] = 0.077601
$data[10,7] = 0.232803
$data[10,8] = 0.02354357842003994
$data[10,9] = 0.02354357842003994
$data[10,10] = 1
$data[10,11] = 0.3333333333333333
$data[10,12] = 0.026642
$data[10,13] = 0.079926
$data[10,14] = 0.001297498167494471
$data[10,15] = 0.001297498167494471
$data[10,16] = 0.002067445842
$data[10,17] = 0.018607012578
$data[10,18] = 0.0000305477498562642
$data[10,19] = 0.0000305477498562642

# Row 13: Inflammatory-Mac -> MuSCs (Efna4/Epha4)
$data[11,0] = "Inflammatory-Mac"
$data[11,1] = "Efna4"
$data[11,2] = "Epha4"
$data[11,3] = "MuSCs"
$data[11,4] = 1
$data[11,5] = 0.3333333333333333
$data[11,6] = 0.077601
$data[11,7] = 0.232803
$data[11,8] = 0.02354357842003994
$data[11,9] = 0.02354357842003994
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 0.391608
$data[11,13] = 1.174824
$data[11,14] = 0.01907179124600912
$data[11,15] = 0.01907179124600912
$data[11,16] = 0.030389172408
$data[11,17] = 0.273502551672
$data[11,18] = 0.0004490182128110469
$data[11,19] = 0.0004490182128110469

# Row 14: MuSCs -> ECs (Efna4/Epha4)
$data[12,0] = "MuSCs"
$data[12,1] = "Efna4"
$data[12,2] = "Epha4"
$data[12,3] = "ECs"
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 0.1837383333333333
$data[12,7] = 0.551215
$data[12,8] = 0.05574487261247628
$data[12,9] = 0.05574487261247628
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 8.540560666666666
$data[12,13] = 25.621682
$data[12,14] = 0.4159358086620884
$data[12,15] = 0.4159358086620884
$data[12,16] = 1.569228382625556
$data[12,17] = 14.12305544363
$data[12,18] = 0.02318628866883542
$data[12,19] = 0.02318628866883542

# Row 15: MuSCs -> FAPs (Efna4/Epha4)
$data[13,0] = "MuSCs"
$data[13,1] = "Efna4"
$data[13,2] = "Epha4"
$data[13,3] = "FAPs"
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 0.1837383333333333
$data[13,7] = 0.551215
$data[13,8] = 0.05574487261247628
$data[13,9] = 0.05574487261247628
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 11.57455166666667
$data[13,13] = 34.723655
$data[13,14] = 0.563694901924408
$data[13,15] = 0.563694901924408
$data[13,16] = 2.126688832313889
$data[13,17] = 19.140199490825
$data[13,18] = 0.03142310050007843
$data[13,19] = 0.03142310050007843

# Row 16: MuSCs -> Inflammatory-Mac (Efna4/Epha4)
$data[14,0] = "MuSCs"
$data[14,1] = "Efna4"
$data[14,2] = "Epha4"
$data[14,3] = "Inflammatory-Mac"
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 0.1837383333333333
$data[14,7] = 0.551215
$data[14,8] = 0.05574487261247628
$data[14,9] = 0.05574487261247628
$data[14,10] = 1
$data[14,11] = 0.3333333333333333
$data[14,12] = 0.026642
$data[14,13] = 0.079926
$data[14,14] = 0.001297498167494471
$data[14,15] = 0.001297498167494471
$data[14,16] = 0.004895156676666667
$data[14,17] = 0.04405641009
$data[14,18] = 0.0000723288700619007
$data[14,19] = 0.0000723288700619007

# Row 17: MuSCs -> MuSCs (Efna4/Epha4)
$data[15,0] = "MuSCs"
$data[15,1] = "Efna4"
$data[15,2] = "Epha4"
$data[15,3] = "MuSCs"
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 0.1837383333333333
$data[15,7] = 0.551215
$data[15,8] = 0.05574487261247628
$data[15,9] = 0.05574487261247628
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 0.391608
$data[15,13] = 1.174824
$data[15,14] = 0.01907179124600912
$data[15,15] = 0.01907179124600912
$data[15,16] = 0.07195340124
$data[15,17] = 0.6475806111600001
$data[15,18] = 0.001063154573500519
$data[15,19] = 0.001063154573500519

# Row 18: Resolving-Mac -> ECs (Efna4/Epha4)
$data[16,0] = "Resolving-Mac"
$data[16,1] = "Efna4"
$data[16,2] = "Epha4"
$data[16,3] = "ECs"
$data[16,4] = 1
$data[16,5] = 0.3333333333333333
$data[16,6] = 0.03989166666666667
$data[16,7] = 0.119675
$data[16,8] = 0.01210284123236505
$data[16,9] = 0.01210284123236505
$data[16,10] = 3
$data[16,11] = 1
$data[16,12] = 8.540560666666666
$data[16,13] = 25.621682
$data[16,14] = 0.4159358086620884
$data[16,15] = 0.4159358086620884
$data[16,16] = 0.3406971992611111
$data[16,17] = 3.06627479335
$data[16,18] = 0.005034005055092622
$data[16,19] = 0.005034005055092623

# Row 19: Resolving-Mac -> FAPs (Efna4/Epha4)
$data[17,0] = "Resolving-Mac"
$data[17,1] = "Efna4"
$data[17,2] = "Epha4"
$data[17,3] = "FAPs"
$data[17,4] = 1
$data[17,5] = 0.3333333333333333
$data[17,6] = 0.03989166666666667
$data[17,7] = 0.119675
$data[17,8] = 0.01210284123236505
$data[17,9] = 0.01210284123236505
$data[17,10] = 3
$data[17,11] = 1
$data[17,12] = 11.57455166666667
$data[17,13] = 34.723655
$data[17,14] = 0.563694901924408
$data[17,15] = 0.563694901924408
$data[17,16] = 0.4617281569027777
$data[17,17] = 4.155553412125
$data[17,18] = 0.006822309901484695
$data[17,19] = 0.006822309901484696

# Row 20: Resolving-Mac -> Inflammatory-Mac (Efna4/Epha4)
$data[18,0] = "Resolving-Mac"
$data[18,1] = "Efna4"
$data[18,2] = "Epha4"
$data[18,3] = "Inflammatory-Mac"
$data[18,4] = 1
$data[18,5] = 0.3333333333333333
$data[18,6] = 0.03989166666666667
$data[18,7] = 0.119675
$data[18,8] = 0.01210284123236505
$data[18,9] = 0.01210284123236505
$data[18,10] = 1
$data[18,11] = 0.3333333333333333
$data[18,12] = 0.026642
$data[18,13] = 0.079926
$data[18,14] = 0.001297498167494471
$data[18,15] = 0.001297498167494471
$data[18,16] = 0.001062793783333333
$data[18,17] = 0.00956514405
$data[18,18] = 0.00001570341432047017
$data[18,19] = 0.00001570341432047018

# Row 21: Resolving-Mac -> MuSCs (Efna4/Epha4)
$data[19,0] = "Resolving-Mac"
$data[19,1] = "Efna4"
$data[19,2] = "Epha4"
$data[19,3] = "MuSCs"
$data[19,4] = 1
$data[19,5] = 0.3333333333333333
$data[19,6] = 0.03989166666666667
$data[19,7] = 0.119675
$data[19,8] = 0.01210284123236505
$data[19,9] = 0.01210284123236505
$data[19,10] = 3
$data[19,11] = 1
$data[19,12] = 0.391608
$data[19,13] = 1.174824
$data[19,14] = 0.01907179124600912
$data[19,15] = 0.01907179124600912
$data[19,16] = 0.0156218958
$data[19,17] = 0.1405970622
$data[19,18] = 0.0002308228614672579
$data[19,19] = 0.0002308228614672579

# Write the whole A2:T21 block in one shot.
$ws.Range("A2:T21").Value2 = $data
